$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing job row (row 2): Java Backend Engineer -> Business Analyst
$ws.Range("B2").Value = "Business Analyst"
$ws.Range("C2").Value = "Remote"
$ws.Range("D2").Value = "Lead requirements gathering, create user stories, drive business process consulting efforts."

# Reset row 2 height back to default (it was auto-tall for the long description before)
$ws.Rows.Item(2).AutoFit()

# Add a new job row (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "QA Analyst"
$ws.Range("C3").Value = "Remote"
$ws.Range("D3").Value = "Condusct functional and regression testing across ServiceNow module, write test cases."

# Update the selection/view state
$ws.Range("D4").Select()
